$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Imperial" to "Manchester"
$ws.Name = "Manchester"

# Update Karen Kirkby's email address to include the ICS alias
$ws.Range("E2").Value = "Karen.kirkby@manchester.ac.uk; karen.kirkby@ics.manchester.ac.uk "

# Remove Michael Merchant's row entirely (author list refresh) -
# this shifts William Bertsche and Stewart Boogert up one row,
# and shrinks the table / used range from A1:P5 to A1:P4.
$ws.Rows.Item(3).Select()
$ws.Rows.Item(3).Delete()
